# Add data for 2021-11-30
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-11-22"

# Row 12 (October) - 2021 columns (T,U,V) updated
$ws.Range("T12").Value = 6
$ws.Range("U12").Value = 189
$ws.Range("V12").Value = 0.0308

# Row 13 (November) - label and all year columns updated
$ws.Range("A13").Value = "November (through 11-22)"

$ws.Range("C13").Value = 21
$ws.Range("D13").Value = 0.0455

$ws.Range("F13").Value = 47
$ws.Range("G13").Value = 0.1132

$ws.Range("I13").Value = 86
$ws.Range("J13").Value = 0.0227

$ws.Range("L13").Value = 37
$ws.Range("M13").Value = 0.1395

$ws.Range("O13").Value = 31
$ws.Range("P13").Value = 0.1389

$ws.Range("R13").Value = 145
$ws.Range("S13").Value = 0.0461

$ws.Range("T13").Value = 3
$ws.Range("U13").Value = 148
$ws.Range("V13").Value = 0.0199

# Row 14 (Total) - recomputed totals across all years
$ws.Range("C14").Value = 247
$ws.Range("D14").Value = 0.1179

$ws.Range("F14").Value = 481
$ws.Range("G14").Value = 0.1076

$ws.Range("I14").Value = 735
$ws.Range("J14").Value = 0.0789

$ws.Range("L14").Value = 586
$ws.Range("M14").Value = 0.1094

$ws.Range("O14").Value = 465
$ws.Range("P14").Value = 0.1023

$ws.Range("R14").Value = 1148
$ws.Range("S14").Value = 0.0505

$ws.Range("T14").Value = 95
$ws.Range("V14").Value = 0.0596
